# Update SA templates with new field
# Inserts a new "Ethnic or Racial Identity" triplet of columns (MODS
# description open-tag, header label, MODS description close-tag) right
# after the existing "Gender Identity, Sexuality" triplet (cols G:I) and
# before "Pronouns" (old cols J:L), shifting everything from old column J
# onward three columns to the right (new cols M: ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns at J:L (old J:L / "Pronouns" triplet and
# everything after slides right to M: onward).
[void]$ws.Range("J1:L1").EntireColumn.Insert()

# Populate the new "Ethnic or Racial Identity" triplet.
$ws.Range("J1").Value = "<mods:description>Ethnic or Racial Identity: "
$ws.Range("K1").Value = "Ethnic or Racial Identity"
$ws.Range("L1").Value = "</mods:description>"

# The header cell (K1) follows the same bold style used by the other
# header cells in the row (e.g. H1 "Gender Identity, Sexuality", N1
# "Pronouns", ...).
$ws.Range("K1").Font.Bold = $true

# Move the active selection to the newly added header cell, matching the
# author's final cursor position.
[void]$ws.Range("J1").Select()
